$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 10.85663951734018
$ws.Range("F2").Value = 2.071808136469829
$ws.Range("G2").Value = 0.1139981867348429
$ws.Range("E3").Value = 10.61493107010667
$ws.Range("F3").Value = 2.041915651538766
$ws.Range("G3").Value = 0.1337238230321848
$ws.Range("E4").Value = 10.36011549237856
$ws.Range("F4").Value = 2.010041405418407
$ws.Range("G4").Value = 0.1545191219416414
$ws.Range("E5").Value = 9.410291888777264
$ws.Range("F5").Value = 1.748210824260051
$ws.Range("G5").Value = 0.2320334792829425
$ws.Range("E6").Value = 9.262309077209016
$ws.Range("F6").Value = 1.711184950795076
$ws.Range("G6").Value = 0.2441102401602038
$ws.Range("E7").Value = 9.09996974905933
$ws.Range("F7").Value = 1.655235964879394
$ws.Range("G7").Value = 0.257358625065601
$ws.Range("E8").Value = 9.380996129053294
$ws.Range("F8").Value = 1.724832870259751
$ws.Range("G8").Value = 0.2344242831956045
$ws.Range("E9").Value = 9.216857865241309
$ws.Range("F9").Value = 1.688579285666962
$ws.Range("G9").Value = 0.247819477825705
$ws.Range("E10").Value = 9.222927468215785
$ws.Range("F10").Value = 1.639054637095778
$ws.Range("G10").Value = 0.2473241423001407
$ws.Range("E11").Value = 10.26973513921639
$ws.Range("F11").Value = 1.972879621920606
$ws.Range("G11").Value = 0.1618949914874004
$ws.Range("E12").Value = 9.944924861996554
$ws.Range("F12").Value = 1.927590761758714
$ws.Range("G12").Value = 0.188402502778006
$ws.Range("E13").Value = 9.586077684113834
$ws.Range("F13").Value = 1.872705225493888
$ws.Range("G13").Value = 0.2176877387648291
$ws.Range("E14").Value = 9.33532629015053
$ws.Range("F14").Value = 1.724388503611734
$ws.Range("G14").Value = 0.2381513628333458
$ws.Range("E15").Value = 9.183411083889451
$ws.Range("F15").Value = 1.692241555634136
$ws.Range("G15").Value = 0.2505490433490235
$ws.Range("E16").Value = 9.082524938475112
$ws.Range("F16").Value = 1.639258138699777
$ws.Range("G16").Value = 0.2587822823386456
$ws.Range("E17").Value = 9.354009065871448
$ws.Range("F17").Value = 1.709156999418435
$ws.Range("G17").Value = 0.2366266761990403
$ws.Range("E18").Value = 9.224596112794856
$ws.Range("F18").Value = 1.677541833946543
$ws.Range("G18").Value = 0.2471879655282776
$ws.Range("E19").Value = 9.388466476436088
$ws.Range("F19").Value = 1.637761978492
$ws.Range("G19").Value = 0.2338146340204335
$ws.Range("E20").Value = 9.766545463004377
$ws.Range("F20").Value = 1.858081106526503
$ws.Range("G20").Value = 0.2029599052508237
$ws.Range("E21").Value = 9.495259572937739
$ws.Range("F21").Value = 1.807212558111466
$ws.Range("G21").Value = 0.2250993333978476
$ws.Range("E22").Value = 9.186579270091167
$ws.Range("F22").Value = 1.742970294113842
$ws.Range("G22").Value = 0.2502904901645875
$ws.Range("E23").Value = 9.304297628528722
$ws.Range("F23").Value = 1.710139452625556
$ws.Range("G23").Value = 0.2406835875071229
$ws.Range("E24").Value = 9.185720605541281
$ws.Range("F24").Value = 1.681482706777054
$ws.Range("G24").Value = 0.2503605651032434
$ws.Range("E25").Value = 9.110261612119324
$ws.Range("F25").Value = 1.619785910127123
$ws.Range("G25").Value = 0.2565187142148754
$ws.Range("E26").Value = 9.335672136843892
$ws.Range("F26").Value = 1.701513438402729
$ws.Range("G26").Value = 0.2381231385566771
$ws.Range("E27").Value = 9.280201746505476
$ws.Range("F27").Value = 1.66286684599664
$ws.Range("G27").Value = 0.2426500334898528
$ws.Range("E28").Value = 9.533565202426065
$ws.Range("F28").Value = 1.635812450048423
$ws.Range("G28").Value = 0.2219732410991476
